$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update F4 and F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2218
$ws1.Range("F6").Value = 369

# Sheet "全部类型" (all types): update F4 and F7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2218
$ws4.Range("F7").Value = 369
